# Converts the "m:Sequence{...}" field code (fldChar/instrText runs) in the
# template into plain literal-text runs: "{m:Sequence{...}}" plus a trailing
# "}" run, matching the TokenIteratorFieldRewriterSplit parser's expected
# input shape. Equivalent to Word's "convert field to text" but emitting an
# extra closing brace, and done by rewriting the paragraph's run content
# directly (Word's Find/Replace does not see hidden field-code text).

$d = $word.ActiveDocument

# Locate the paragraph that holds the MERGEFIELD-style field (rather than
# assuming a fixed paragraph index).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Fields.Count -gt 0) {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the paragraph containing the m:Sequence field"
}

$newParagraphXml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" w:rsidP="00F5495F" w:rsidR="00C52979" w:rsidRDefault="00C52979"><w:r><w:t>{</w:t></w:r><w:r><w:t>m</w:t></w:r><w:r><w:t>:</w:t></w:r><w:r><w:t>Sequence{</w:t></w:r><w:r><w:t>'dh1.gif'.asImage()</w:t></w:r><w:r><w:t>.setWidth(50)</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t>'dh1.gif'.asImage()</w:t></w:r><w:r><w:t>.setWidth(50)</w:t></w:r><w:r><w:t>,</w:t></w:r><w:r><w:t>'dh1.gif'.asImage()</w:t></w:r><w:r><w:t>.setWidth(50)</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>}</w:t></w:r><w:r><w:t xml:space="preserve">}</w:t></w:r></w:p>
'@

$target.Range.InsertXML($newParagraphXml)
